$d = $word.ActiveDocument

# Locate the "NO FUNCIONALES" heading paragraph; the new R8 table + spacer paragraph
# must be inserted immediately before it (right after the spacer paragraph that
# currently follows the R7 table).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "FUNCIONALES") {
        if ($p.Range.Text -match "NO") {
            $target = $p
            break
        }
    }
}

$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2263"/><w:gridCol w:w="7087"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2263" w:type="dxa"/><w:shd w:val="solid" w:color="70AD47" w:themeColor="accent6" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Nombre </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="7087" w:type="dxa"/><w:shd w:val="solid" w:color="595959" w:themeColor="text1" w:themeTint="A6" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t>8</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Respaldo </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2263" w:type="dxa"/><w:shd w:val="solid" w:color="70AD47" w:themeColor="accent6" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t>Resumen</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="7087" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr><w:t>La aplicación debe porporcionar un respaldo de todos los datos y operaciones realizados por la aplicación</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr><w:t>.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:gridAfter w:val="1"/><w:wAfter w:w="7087" w:type="dxa"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2263" w:type="dxa"/><w:shd w:val="solid" w:color="70AD47" w:themeColor="accent6" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t>Entradas</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="125"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2263" w:type="dxa"/><w:shd w:val="solid" w:color="70AD47" w:themeColor="accent6" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="es-CO"/></w:rPr><w:t>Resultado</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="7087" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Todos los datos de interés son alojados en la base de datos. </w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng = $target.Range
$rng.Collapse(1)
$rng.InsertXML($insertXml)

Write-Host "Inserted. Tables now:" $d.Tables.Count
